$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Move Robot21 to location (6, 12) and remove the toolkit."
$ws.Range("B1").Value = "['Robot22']"
$ws.Range("E1").Value = "(6, 12)"

# Row 2
$ws.Range("A2").Value = "Move Robot41 to location (1, 11) and remove the liquid spill."
$ws.Range("B2").Value = "['Robot6']"
$ws.Range("E2").Value = "(1, 11)"

# Row 3
$ws.Range("A3").Value = "Move Robot9 to location (1, 5) and remove the large debris."
$ws.Range("B3").Value = "['Robot2', 'Robot39']"
$ws.Range("E3").Value = "(1, 5)"

# Row 4
$ws.Range("A4").Value = "Move Robot42 to location (1, 11) and remove the dust."
$ws.Range("E4").Value = "(1, 11)"

# Row 5
$ws.Range("A5").Value = "Move Robot32 to location (3, 9) and remove the grass."
$ws.Range("B5").Value = "['Robot21']"
$ws.Range("E5").Value = "(3, 9)"

# Row 6
$ws.Range("A6").Value = "Move Robot14 to location (11, 12) and remove the small debris."
$ws.Range("B6").Value = "['Robot8', 'Robot10']"
$ws.Range("E6").Value = "(11, 12)"

# Row 7
$ws.Range("A7").Value = "Move Robot39 to location (6, 4) and remove the vehicle."
$ws.Range("E7").Value = "(6, 4)"

# Row 8
$ws.Range("A8").Value = "Move Robot15 to location (11, 2) and remove the construction materials."
$ws.Range("E8").Value = "(11, 2)"

# Row 9
$ws.Range("A9").Value = "Move Robot2 to location (1, 10) and remove the tree branches."
$ws.Range("B9").Value = "['Robot14']"
$ws.Range("E9").Value = "(1, 10)"

# Row 10
$ws.Range("A10").Value = "Move Robot26 to location (1, 3) and remove the screws."
$ws.Range("E10").Value = "(1, 3)"
